{"js": "// Replace each math-problem answer cell's text with its updated value.\n// The document body contains a single 20-row x 5-column table where every\n// cell holds one line of text like \"90-66=24\". The commit replaces every\n// cell's text with a new expression (values below, in row-major reading\n// order matching the table's natural cell order).\nconst newValues = [\"83-55=28\", \"19+73=92\", \"81-72=9\", \"69+14=83\", \"55+18=73\", \"63-5=58\", \"27+18=45\", \"29+58=87\", \"84+8=92\", \"8+59=67\", \"60-51=9\", \"57-9=48\", \"57-48=9\", \"58+24=82\", \"34+8=42\", \"90-75=15\", \"93-36=57\", \"25+47=72\", \"7+39=46\", \"30-25=5\", \"94-45=49\", \"29+4=33\", \"47+8=55\", \"25+37=62\", \"57+25=82\", \"46+19=65\", \"55-39=16\", \"39+49=88\", \"69+6=75\", \"4+8=12\", \"29+23=52\", \"8+25=33\", \"69+5=74\", \"67-28=39\", \"14+8=22\", \"57+28=85\", \"64-7=57\", \"79+16=95\", \"31-16=15\", \"12+49=61\", \"18+15=33\", \"83-24=59\", \"90-27=63\", \"93-78=15\", \"85-16=69\", \"31-3=28\", \"76-58=18\", \"8+19=27\", \"94-55=39\", \"87+4=91\", \"84-47=37\", \"57+5=62\", \"29+18=47\", \"74-29=45\", \"14+48=62\", \"48+45=93\", \"76-49=27\", \"90-26=64\", \"71-64=7\", \"8+53=61\", \"46-27=19\", \"26+47=73\", \"15-8=7\", \"50-25=25\", \"37+46=83\", \"16+26=42\", \"25+67=92\", \"15-8=7\", \"27-18=9\", \"81-46=35\", \"82-35=47\", \"45+18=63\", \"41-2=39\", \"62-35=27\", \"83+8=91\", \"32-19=13\", \"97-19=78\", \"71-62=9\", \"77+5=82\", \"35-28=7\", \"50-44=6\", \"67+27=94\", \"35+56=91\", \"46+8=54\", \"18+49=67\", \"61-56=5\", \"61-57=4\", \"82-75=7\", \"19+18=37\", \"61-19=42\", \"28+56=84\", \"30-28=2\", \"10-6=4\", \"10-7=3\", \"68+7=75\", \"37+35=72\", \"91-43=48\", \"80-36=44\", \"63-8=55\", \"9+77=86\"];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst cols = table.values && table.values[0] ? table.values[0].length : 5;\nfor (let i = 0; i < newValues.length; i++) {\n  const row = Math.floor(i / cols);\n  const col = i % cols;\n  table.getCell(row, col).value = newValues[i];\n}\nawait context.sync();\n", "ps1": "# Replace each math-problem answer cell's text with its updated value.\n# The document body contains a single 20-row x 5-column table where every\n# cell holds one line of text like \"90-66=24\". The commit replaces every\n# cell's text with a new expression (values below, in row-major reading\n# order matching the table's natural cell order).\n$newValues = @(\n    \"83-55=28\",\n    \"19+73=92\",\n    \"81-72=9\",\n    \"69+14=83\",\n    \"55+18=73\",\n    \"63-5=58\",\n    \"27+18=45\",\n    \"29+58=87\",\n    \"84+8=92\",\n    \"8+59=67\",\n    \"60-51=9\",\n    \"57-9=48\",\n    \"57-48=9\",\n    \"58+24=82\",\n    \"34+8=42\",\n    \"90-75=15\",\n    \"93-36=57\",\n    \"25+47=72\",\n    \"7+39=46\",\n    \"30-25=5\",\n    \"94-45=49\",\n    \"29+4=33\",\n    \"47+8=55\",\n    \"25+37=62\",\n    \"57+25=82\",\n    \"46+19=65\",\n    \"55-39=16\",\n    \"39+49=88\",\n    \"69+6=75\",\n    \"4+8=12\",\n    \"29+23=52\",\n    \"8+25=33\",\n    \"69+5=74\",\n    \"67-28=39\",\n    \"14+8=22\",\n    \"57+28=85\",\n    \"64-7=57\",\n    \"79+16=95\",\n    \"31-16=15\",\n    \"12+49=61\",\n    \"18+15=33\",\n    \"83-24=59\",\n    \"90-27=63\",\n    \"93-78=15\",\n    \"85-16=69\",\n    \"31-3=28\",\n    \"76-58=18\",\n    \"8+19=27\",\n    \"94-55=39\",\n    \"87+4=91\",\n    \"84-47=37\",\n    \"57+5=62\",\n    \"29+18=47\",\n    \"74-29=45\",\n    \"14+48=62\",\n    \"48+45=93\",\n    \"76-49=27\",\n    \"90-26=64\",\n    \"71-64=7\",\n    \"8+53=61\",\n    \"46-27=19\",\n    \"26+47=73\",\n    \"15-8=7\",\n    \"50-25=25\",\n    \"37+46=83\",\n    \"16+26=42\",\n    \"25+67=92\",\n    \"15-8=7\",\n    \"27-18=9\",\n    \"81-46=35\",\n    \"82-35=47\",\n    \"45+18=63\",\n    \"41-2=39\",\n    \"62-35=27\",\n    \"83+8=91\",\n    \"32-19=13\",\n    \"97-19=78\",\n    \"71-62=9\",\n    \"77+5=82\",\n    \"35-28=7\",\n    \"50-44=6\",\n    \"67+27=94\",\n    \"35+56=91\",\n    \"46+8=54\",\n    \"18+49=67\",\n    \"61-56=5\",\n    \"61-57=4\",\n    \"82-75=7\",\n    \"19+18=37\",\n    \"61-19=42\",\n    \"28+56=84\",\n    \"30-28=2\",\n    \"10-6=4\",\n    \"10-7=3\",\n    \"68+7=75\",\n    \"37+35=72\",\n    \"91-43=48\",\n    \"80-36=44\",\n    \"63-8=55\",\n    \"9+77=86\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$cols = $t.Columns.Count\n\nfor ($i = 0; $i -lt $newValues.Count; $i++) {\n    $row = [Math]::Floor($i / $cols) + 1\n    $col = ($i % $cols) + 1\n    $t.Cell($row, $col).Range.Text = $newValues[$i]\n}\n"}
